$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20: mark user as validated (isValidated -> TRUE) and set their
# userType-specific attribute (column H) to "Computer Science"
$ws.Range("F20").Value = $true
$ws.Range("H20").Value = "Computer Science"
